$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet0")

$newIds = @(
    "FT231680B2S7WCZ1",
    "FT231680N8731K99",
    "FT231680NBJ05K5C",
    "FT231680NBJ05LMJ",
    "FT231680G1NNFD6X"
)

$startRow = 18
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}
